# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to match the regenerated gh-pages data output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$wsExhibit = $wb.Worksheets.Item("展览")

$wsExhibit.Range("F2").Value  = 6939
$wsExhibit.Range("F4").Value  = 0
$wsExhibit.Range("F6").Value  = 0
$wsExhibit.Range("F9").Value  = 0
$wsExhibit.Range("F10").Value = 0
$wsExhibit.Range("F11").Value = 0
$wsExhibit.Range("F16").Value = 407
$wsExhibit.Range("F17").Value = 48
$wsExhibit.Range("F18").Value = 27
$wsExhibit.Range("F19").Value = 0
$wsExhibit.Range("F24").Value = 0
$wsExhibit.Range("F25").Value = 220

# ---- Sheet: 全部类型 ----
$wsAll = $wb.Worksheets.Item("全部类型")

$wsAll.Range("F2").Value  = 6939
$wsAll.Range("F3").Value  = 98
$wsAll.Range("F4").Value  = 51
$wsAll.Range("F8").Value  = 68
$wsAll.Range("F9").Value  = 0
$wsAll.Range("F13").Value = 405
$wsAll.Range("F14").Value = 0
$wsAll.Range("F17").Value = 48
$wsAll.Range("F18").Value = 0
$wsAll.Range("F19").Value = 14
$wsAll.Range("F21").Value = 0
$wsAll.Range("F24").Value = 154
$wsAll.Range("F25").Value = 568

$wb.Save()
